# The edit adds a new "bid_made" column (Q) to the "Plays" sheet:
#   - Q2 gets the header "bid_made"
#   - Q3 gets a boolean TRUE value
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plays")

$ws.Range("Q2").Value = "bid_made"
$ws.Range("Q3").Value = $true
